# ============================================================================
# Edit: creepy-carnival review document
#  1. Insert a new "Meta description" paragraph right after the title
#     (Heading1) paragraph, with "Meta description" in bold followed by a
#     plain-text blurb.
#  2. Near the bottom of the document, delete the paragraph that duplicated
#     the page title in bold, and replace the text of the following italic
#     paragraph with the new image-generation prompt text.
# ============================================================================

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: insert the "Meta description" paragraph after the Heading1 title.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)

# Split the title paragraph's own range so a brand-new (empty) paragraph is
# created right after it. The new paragraph temporarily inherits the
# Heading1 style/run content, but it is about to be fully replaced below.
$titlePara.Range.InsertParagraphAfter() | Out-Null
$newPara = $d.Paragraphs(2)
$newParaRange = $newPara.Range.Duplicate

$metaXml = @'
<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Explore the eerie circus show of NoLimit City&apos;s Creepy Carnival. Play for free and enjoy two unique features - the Free Spin and Star Spin mode.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

# Replacing the *entire* range of the freshly-created paragraph (not just a
# collapsed insertion point) swaps in our new paragraph cleanly, without
# disturbing the title paragraph before it or the "Game Overview" paragraph
# that follows it.
$newParaRange.InsertXML($metaXml) | Out-Null

# ---------------------------------------------------------------------------
# Step 2: near the end of the document, drop the paragraph that duplicated
# the page title in bold, and update the following italic paragraph's text.
# ---------------------------------------------------------------------------
$boldTitleText = "Play Creepy Carnival Free - Review of NoLimit City's Spooky Slot"

$boldTitlePara = $null
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.TrimEnd([char]13) -eq $boldTitleText -and $i -ne 1) {
        $boldTitlePara = $p
        break
    }
}
$boldTitlePara.Range.Delete()

# The paragraph that used to follow the bold duplicated title is now the
# last paragraph in the document; update its (italic) text in place so its
# run formatting is preserved.
$italicPara = $d.Paragraphs($d.Paragraphs.Count)
$italicRange = $italicPara.Range.Duplicate
$italicRange.MoveEnd(1, -1) | Out-Null
$italicRange.Text = "Prompt: Create a feature image for Creepy Carnival that fits the game using the following criteria: - Cartoon style image - Happy Maya warrior with glasses For the feature image of Creepy Carnival, we want to create a playful and cartoonish depiction of the game that still highlights its eerie and mysterious atmosphere. We want to incorporate a happy Maya warrior with glasses to symbolize the adventure and thrill of exploring the dark and twisted world of the game. The image will feature the Maya warrior standing in front of a creepy carnival backdrop with strange and terrifying creatures surrounding him. He will be holding a glowing crystal ball in one hand and an ancient Mayan artifact in the other, both of which represent the mystical and supernatural elements of the game. The Maya warrior will be wearing glasses, giving him a modern and cool edge while also emphasizing his intelligence and courage. His happy expression will add a touch of joy and excitement to the otherwise spooky and ominous image. Overall, the feature image will be colorful and attention-grabbing, conveying the sense of adventure and danger that players will experience while playing Creepy Carnival. It will capture the essence of the game in a fun and engaging way, inviting players to enter the creepy and fascinating world of the carnival."

Write-Output "Edits applied successfully"
